# Append new Order_Items rows (4-8) to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(3, 3, "BP-XL-004",  7, 31500, 0),
    @(4, 3, "BP-M-002",   4, 10000, 0),
    @(5, 4, "BP-M-002",   6, 15000, 0),
    @(6, 4, "BP-XXL-005", 8, 44000, 0),
    @(7, 5, "CL-G-008",   3,  7500, 0)
)

$startRow = 4
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
